$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 31: Nutr_No=221 (ethanol), unit=g, Tagname=ALC, NutrDesc=Alcohol, ethyl, nutCode=ethanol_g
$ws.Range("A31").Value = 221
$ws.Range("B31").Value = "g"
$ws.Range("C31").Value = "ALC"
$ws.Range("D31").Value = "Alcohol, ethyl"
$ws.Range("E31").Value = "ethanol_g"

# NutrDesc cell uses a distinct font (Lucida Sans, 11pt, black)
$ws.Range("D31").Font.Name = "Lucida Sans"
$ws.Range("D31").Font.Size = 11
$ws.Range("D31").Font.Color = 0

# Selection moves to the newly entered cell
$ws.Range("C31").Select()
